$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.606.84'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '2.341.45'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '305.03'
$c.Style = "Normal"

$ws.Range("E5").Value = '  -1.80%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '102.04'
$c.Style = "Normal"

$ws.Range("E6").Value = '  -2.65%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.512'
$c.Style = "Normal"

$ws.Range("E7").Value = '  -2.42%  '
$ws.Range("E8").Value = '  +0.03%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.511'
$c.Style = "Normal"

$ws.Range("E9").Value = '  -1.40%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '35.28'
$c.Style = "Normal"

$ws.Range("E10").Value = '  -2.83%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0798'
$c.Style = "Normal"

$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("D14").Value = '2.704.95'
$ws.Range("E14").Value = '  -1.07%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.58'
$c.Style = "Normal"

$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").Value = '2.374.59'
$ws.Range("E16").Value = '  -0.16%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.808'
$c.Style = "Normal"

$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = '43.509.80'
$ws.Range("E18").Value = '  +0.37%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.88'
$c.Style = "Normal"

$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("E20").Value = '  -1.92%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.12'
$c.Style = "Normal"

$ws.Range("E21").Value = '  -2.64%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '68.33'
$c.Style = "Normal"

$ws.Range("E22").Value = '  -0.17%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '238.45'
$c.Style = "Normal"

$ws.Range("E23").Value = '  -1.47%  '
$ws.Range("E24").Value = '  -3.64%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.55'
$c.Style = "Normal"

$ws.Range("E25").Value = '  -2.89%  '
$ws.Range("E26").Value = '  +0.03%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '25.05'
$c.Style = "Normal"

$ws.Range("E27").Value = '  -3.45%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '34.65'
$c.Style = "Normal"

$ws.Range("E28").Value = '  -5.94%  '
$ws.Range("E29").Value = '  -5.88%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '9.29'
$c.Style = "Normal"

$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '166.16'
$c.Style = "Normal"

$ws.Range("E31").Value = '  +2.66%  '
$ws.Range("E32").Value = '  -0.01%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.07'
$c.Style = "Normal"

$ws.Range("E33").Value = '  -4.23%  '
$ws.Range("E34").Value = '  -4.83%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.50'
$c.Style = "Normal"

$ws.Range("E35").Value = '  -3.96%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '16.97'
$c.Style = "Normal"

$ws.Range("E36").Value = '  -7.67%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.0708'
$c.Style = "Normal"

$ws.Range("E37").Value = '  -4.44%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.91'
$c.Style = "Normal"

$ws.Range("E38").Value = '  -6.65%  '
$ws.Range("E39").Value = '  -6.27%  '
$ws.Range("E40").Value = '  -2.74%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.111'
$c.Style = "Normal"

$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("E42").Value = '  -1.93%  '
$ws.Range("D43").Value = '1.992.51'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("E44").Value = '  -1.99%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '18.48'
$c.Style = "Normal"

$ws.Range("E45").Value = '  -8.66%  '
$ws.Range("E46").Value = '  -6.94%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.00'
$c.Style = "Normal"

$ws.Range("E47").Value = '  -3.82%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '56.62'
$c.Style = "Normal"

$ws.Range("E48").Value = '  -2.83%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '4.91'
$c.Style = "Normal"

$ws.Range("E49").Value = '  +4.22%  '
$ws.Range("D50").Value = '2.566.30'
$ws.Range("E51").Value = '  -1.18%  '
